$d = $word.ActiveDocument

# Update the date line (unique text in the document body)
$d.Content.Find.Execute("2025-12-13 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-14 Sunday", 2) | Out-Null

# Update the division problems in the table, cell by cell (positional,
# since several cells share identical source text e.g. 56÷8=)
$tbl = $d.Tables.Item(1)
$tbl.Cell(1, 1).Range.Text = "52÷6="
$tbl.Cell(1, 2).Range.Text = "79÷8="
$tbl.Cell(1, 3).Range.Text = "11÷9="
$tbl.Cell(1, 4).Range.Text = "95÷9="
$tbl.Cell(1, 5).Range.Text = "82÷6="
$tbl.Cell(5, 1).Range.Text = "47÷8="
$tbl.Cell(5, 2).Range.Text = "78÷3="
$tbl.Cell(5, 3).Range.Text = "50÷8="
$tbl.Cell(5, 4).Range.Text = "78÷3="
$tbl.Cell(5, 5).Range.Text = "20÷3="
$tbl.Cell(9, 1).Range.Text = "90÷3="
$tbl.Cell(9, 2).Range.Text = "21÷2="
$tbl.Cell(9, 3).Range.Text = "94÷3="
$tbl.Cell(9, 4).Range.Text = "17÷5="
$tbl.Cell(9, 5).Range.Text = "14÷4="
$tbl.Cell(13, 1).Range.Text = "92÷2="
$tbl.Cell(13, 2).Range.Text = "60÷3="
$tbl.Cell(13, 3).Range.Text = "89÷4="
$tbl.Cell(13, 4).Range.Text = "57÷9="
$tbl.Cell(13, 5).Range.Text = "27÷8="
$tbl.Cell(17, 1).Range.Text = "17÷6="
$tbl.Cell(17, 2).Range.Text = "18÷8="
$tbl.Cell(17, 3).Range.Text = "44÷5="
$tbl.Cell(17, 4).Range.Text = "42÷8="
$tbl.Cell(17, 5).Range.Text = "51÷9="
